$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Finish French first 8 moves: fill in the three remaining French chapters ---

$ws.Range("D13").Value = "E:\Chess\Database\Openings\French-Defense\French-MacCutcheon.pgn"
$ws.Range("E13").Value = "Done*"
$ws.Range("E13").Font.Bold = $true
$ws.Range("F13").Value = "All lines followed at least to move 8"

$ws.Range("D14").Value = "E:\Chess\Database\Openings\French-Defense\French-Classical.pgn"
$ws.Range("E14").Value = "Done*"
$ws.Range("E14").Font.Bold = $true
$ws.Range("F14").Value = "All lines followed at least to move 8"

$ws.Range("D15").Value = "E:\Chess\Database\Openings\French-Defense\French-Winawer.pgn"
$ws.Range("E15").Value = "Done*"
$ws.Range("E15").Font.Bold = $true
$ws.Range("F15").Value = "All lines followed at least to move 8"

# --- Add a new (blank) Sheet2 after Sheet1 ---

$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Sheet2"

# --- Restore Sheet1 as active and move the cursor to where editing left off ---

$ws.Activate()
$null = $ws.Range("D16").Select()
